$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-5 belong to the "Hospitality" department (previously "FACULTY OF HOSPITALITY")
$ws.Range("C2").Value = "Hospitality"
$ws.Range("C3").Value = "Hospitality"
$ws.Range("C4").Value = "Hospitality"
$ws.Range("C5").Value = "Hospitality"

# Rows 6-8 belong to the "Packages" department (previously "FACULTY OF HOSPITALITY")
$ws.Range("C6").Value = "Packages"
$ws.Range("C7").Value = "Packages"
$ws.Range("C8").Value = "Packages"

# Row heights settle to a uniform 42.75 after the department label text shortens
$ws.Rows.Item(2).RowHeight = 42.75
$ws.Rows.Item(3).RowHeight = 42.75
$ws.Rows.Item(4).RowHeight = 42.75
$ws.Rows.Item(5).RowHeight = 42.75
$ws.Rows.Item(6).RowHeight = 42.75
$ws.Rows.Item(7).RowHeight = 42.75
$ws.Rows.Item(8).RowHeight = 42.75

# Selection moves to C8 in the saved view state
$ws.Range("C8").Select() | Out-Null
